$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.131.08"
$ws.Cells.Item(2, 5).Value = "  +0.09%  "
$ws.Cells.Item(3, 4).Value = "1.799.79"
$ws.Cells.Item(3, 5).Value = "  +2.31%  "
$ws.Cells.Item(4, 4).Value = "'1.005"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.19%  "
$ws.Cells.Item(5, 4).Value = "'338.00"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.76%  "
$ws.Cells.Item(6, 5).Value = "  +0.12%  "
$ws.Cells.Item(7, 4).Value = "'0.4745"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +25.50%  "
$ws.Cells.Item(8, 4).Value = "'0.3739"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +11.23%  "
$ws.Cells.Item(9, 4).Value = "'45.58"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.44%  "
$ws.Cells.Item(10, 4).Value = "'0.07685"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +6.38%  "
$ws.Cells.Item(11, 4).Value = "'1.148"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.46%  "
$ws.Cells.Item(12, 4).Value = "'22.63"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.57%  "
$ws.Cells.Item(13, 4).Value = "'1.003"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.26%  "
$ws.Cells.Item(14, 5).Value = "  +3.28%  "
$ws.Cells.Item(15, 5).Value = "  +3.33%  "
$ws.Cells.Item(16, 4).Value = "1.796.31"
$ws.Cells.Item(16, 5).Value = "  +2.15%  "
$ws.Cells.Item(17, 5).Value = "  +3.64%  "
$ws.Cells.Item(18, 4).Value = "'0.06736"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.97%  "
$ws.Cells.Item(19, 4).Value = "'82.62"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +2.01%  "
$ws.Cells.Item(20, 5).Value = "  +0.12%  "
$ws.Cells.Item(21, 4).Value = "'17.51"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.78%  "
$ws.Cells.Item(22, 4).Value = "'6.443"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.61%  "
$ws.Cells.Item(23, 4).Value = "28.145.55"
$ws.Cells.Item(23, 5).Value = "  +0.17%  "
$ws.Cells.Item(24, 4).Value = "'11.99"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.39%  "
$ws.Cells.Item(25, 4).Value = "'2.402"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.30%  "
$ws.Cells.Item(26, 4).Value = "'20.94"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +5.05%  "
$ws.Cells.Item(27, 4).Value = "'2.407"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.83%  "
$ws.Cells.Item(28, 4).Value = "'152.02"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.05%  "
$ws.Cells.Item(29, 4).Value = "2.002.16"
$ws.Cells.Item(30, 4).Value = "'134.34"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.57%  "
$ws.Cells.Item(31, 4).Value = "'1.266"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.15%  "
$ws.Cells.Item(32, 4).Value = "'4.051"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.77%  "
$ws.Cells.Item(33, 4).Value = "'0.09674"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +10.39%  "
$ws.Cells.Item(34, 4).Value = "'5.956"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.14%  "
$ws.Cells.Item(35, 4).Value = "'0.02403"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.26%  "
$ws.Cells.Item(36, 4).Value = "'12.25"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.75%  "
$ws.Cells.Item(37, 4).Value = "'0.2230"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +5.04%  "
$ws.Cells.Item(38, 4).Value = "'0.06381"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.29%  "
$ws.Cells.Item(39, 5).Value = "  +0.63%  "
$ws.Cells.Item(40, 4).Value = "'5.278"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.05%  "
$ws.Cells.Item(41, 4).Value = "'1.237"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.36%  "
$ws.Cells.Item(42, 5).Value = "  +2.59%  "
$ws.Cells.Item(43, 4).Value = "'8.111"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.74%  "
$ws.Cells.Item(44, 5).Value = "  +4.10%  "
$ws.Cells.Item(45, 5).Value = "  +0.10%  "
$ws.Cells.Item(46, 4).Value = "'0.6175"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.66%  "
$ws.Cells.Item(47, 4).Value = "'3.862"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.18%  "
$ws.Cells.Item(48, 4).Value = "'130.37"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.52%  "
$ws.Cells.Item(49, 4).Value = "'2.071"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.33%  "
$ws.Cells.Item(50, 4).Value = "'1.184"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.77%  "
$ws.Cells.Item(51, 4).Value = "'0.07127"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.28%  "
